$d = $word.ActiveDocument

# Helper-style replacements via Find/Replace across the whole document body.
# wdFindContinue = 1, wdReplaceAll = 2

# 1. Title / headline - used both as the H1 heading and later as a bold
#    "recap" line near the end of the document (Replace:=2 handles both).
$d.Content.Find.Execute(
    "Play Fever for Free - Slot Game Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Fever for Free - Exciting Slot Game with Modern Features", 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute(
    "Combines classic mechanics with modern features", $true, $false, $false, $false, $false,
    $true, 1, $false, "Combines classic slot mechanics with modern features", 2)

$d.Content.Find.Execute(
    "Appealing to a wider range of users", $true, $false, $false, $false, $false,
    $true, 1, $false, "Appealing design and audio with disco-style setting", 2)

$d.Content.Find.Execute(
    "Simple yet enjoyable gameplay", $true, $false, $false, $false, $false,
    $true, 1, $false, "Special symbols increase chances of winning", 2)

$d.Content.Find.Execute(
    "Well-executed design and audio", $true, $false, $false, $false, $false,
    $true, 1, $false, "Easy to understand gameplay, ideal for novice players", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute(
    "May seem a bit bare", $true, $false, $false, $false, $false,
    $true, 1, $false, "Game may seem a bit bare", 2)

$d.Content.Find.Execute(
    "Follows classic slot archetypes", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited variety of symbols", 2)

# 4. Meta description paragraph (italic run near the end)
$d.Content.Find.Execute(
    "Read our review of Fever by Cristaltec, a slot game that combines classic mechanics with modern features. Play for free and trigger free spins!",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Read our review of Fever, a slot game that combines classic mechanics with modern features. Play for free and enjoy the immersive disco-style setting.",
    2)
